$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 724.25
$ws.Range("I2").Value = 816.6667
$ws.Range("J2").Value = 447
$ws.Range("K2").Value = 816.6667
$ws.Range("L2").Value = 447
$ws.Range("M2").Value = -703.6667
$ws.Range("N2").Value = -673

$ws.Range("H40").Value = 2054.2222
$ws.Range("I40").Value = 1213.4286
$ws.Range("J40").Value = 4997
$ws.Range("K40").Value = 1213.4286
$ws.Range("L40").Value = 4997
$ws.Range("M40").Value = -1038.4286
$ws.Range("N40").Value = -5347

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = $null

$ws.Range("H137").Value = 3307.389
$ws.Range("I137").Value = 2727.0625
$ws.Range("J137").Value = 7950
$ws.Range("K137").Value = 8181.1875
$ws.Range("L137").Value = 23850
$ws.Range("M137").Value = -5631.1875
$ws.Range("N137").Value = -28950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1414.7273
$ws.Range("I2").Value = 1007
$ws.Range("J2").Value = 2128.25
$ws.Range("K2").Value = 1007
$ws.Range("L2").Value = 2128.25
$ws.Range("M2").Value = -894
$ws.Range("N2").Value = -2354.25

$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("K14").Value = 500
$ws.Range("M14").Value = -325

$ws.Range("H102").Value = 1775.5555
$ws.Range("I102").Value = 1775.5555
$ws.Range("K102").Value = 1775.5555
$ws.Range("M102").Value = -153.5554999999999

$ws.Range("H110").Value = 7823.857
$ws.Range("I110").Value = 8961.166999999999
$ws.Range("K110").Value = 8961.166999999999
$ws.Range("M110").Value = -6916.166999999999

$ws.Range("H116").Value = 1414.7273
$ws.Range("I116").Value = 1007
$ws.Range("J116").Value = 2128.25
$ws.Range("K116").Value = 1007
$ws.Range("L116").Value = 2128.25
$ws.Range("M116").Value = 1287
$ws.Range("N116").Value = -6716.25

$ws.Range("H132").Value = 3089.7778
$ws.Range("I132").Value = 3089.7778
$ws.Range("K132").Value = 9269.3334
$ws.Range("M132").Value = -6739.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1414.7273
$ws.Range("I3").Value = 1007
$ws.Range("J3").Value = 2128.25
$ws.Range("K3").Value = 1007
$ws.Range("L3").Value = 2128.25
$ws.Range("M3").Value = -893
$ws.Range("N3").Value = -2356.25

$ws.Range("H20").Value = 1881.8823
$ws.Range("I20").Value = 1191.5
$ws.Range("J20").Value = 2868.1428
$ws.Range("K20").Value = 1191.5
$ws.Range("L20").Value = 2868.1428
$ws.Range("M20").Value = -944.5
$ws.Range("N20").Value = -3362.1428

$ws.Range("H105").Value = 4754.2354
$ws.Range("I105").Value = 4280.4287
$ws.Range("K105").Value = 4280.4287
$ws.Range("M105").Value = -2533.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2310
$ws.Range("I31").Value = 2310
$ws.Range("K31").Value = 2310
$ws.Range("M31").Value = -2015

$ws.Range("H34").Value = 2310
$ws.Range("I34").Value = 2310
$ws.Range("K34").Value = 2310
$ws.Range("M34").Value = -2108

$ws.Range("H60").Value = 45724.5
$ws.Range("J60").Value = 24998
$ws.Range("L60").Value = 24998
$ws.Range("N60").Value = -26020

$ws.Range("H86").Value = 7597
$ws.Range("I86").Value = 7492.5
$ws.Range("K86").Value = 7492.5
$ws.Range("M86").Value = -6369.5

$ws.Range("H89").Value = 7597
$ws.Range("I89").Value = 7492.5
$ws.Range("K89").Value = 37462.5
$ws.Range("M89").Value = -31846.5

$ws.Range("H94").Value = 630.4545000000001
$ws.Range("I94").Value = 897.6
$ws.Range("J94").Value = 407.83334
$ws.Range("K94").Value = 897.6
$ws.Range("L94").Value = 407.83334
$ws.Range("M94").Value = -446.6
$ws.Range("N94").Value = -1309.83334

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null

$ws.Range("H107").Value = 1316.375
$ws.Range("I107").Value = 758
$ws.Range("J107").Value = 2991.5
$ws.Range("K107").Value = 758
$ws.Range("L107").Value = 2991.5
$ws.Range("M107").Value = 1162
$ws.Range("N107").Value = -6831.5

$ws.Range("H122").Value = 3651.0833
$ws.Range("I122").Value = 4333.1665
$ws.Range("J122").Value = 2969
$ws.Range("K122").Value = 12999.4995
$ws.Range("L122").Value = 8907
$ws.Range("M122").Value = -10549.4995
$ws.Range("N122").Value = -13807

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 352.5
$ws.Range("I64").Value = 352.5
$ws.Range("K64").Value = 1057.5
$ws.Range("M64").Value = -787.5

$ws.Range("H67").Value = 352.5
$ws.Range("I67").Value = 352.5
$ws.Range("K67").Value = 1057.5
$ws.Range("M67").Value = -121.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1553.091
$ws.Range("I102").Value = 1553.091
$ws.Range("K102").Value = 1553.091
$ws.Range("M102").Value = 68.90900000000011

$ws.Range("H122").Value = 2024.3334
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 5983
$ws.Range("I132").Value = 5983
$ws.Range("K132").Value = 17949
$ws.Range("M132").Value = -15419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9496
$ws.Range("I7").Value = 9496
$ws.Range("K7").Value = 9496
$ws.Range("M7").Value = -9384

$ws.Range("H122").Value = 4492.8667
$ws.Range("I122").Value = 3308.5454
$ws.Range("K122").Value = 9925.636200000001
$ws.Range("M122").Value = -7475.636200000001

$ws.Range("H126").Value = 9496
$ws.Range("I126").Value = 9496
$ws.Range("K126").Value = 28488
$ws.Range("M126").Value = -26018

$ws.Range("H132").Value = 5465.2
$ws.Range("I132").Value = 5021
$ws.Range("K132").Value = 15063
$ws.Range("M132").Value = -12533

$ws.Range("H136").Value = 35716276
$ws.Range("I136").Value = 1741.1818
$ws.Range("J136").Value = 166669570
$ws.Range("K136").Value = 5223.5454
$ws.Range("L136").Value = 500008710
$ws.Range("M136").Value = -2673.5454
$ws.Range("N136").Value = -500013810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null

$ws.Range("H104").Value = 24163
$ws.Range("J104").Value = 24163
$ws.Range("L104").Value = 24163
$ws.Range("N104").Value = -31151

$ws.Range("H122").Value = 5816.923
$ws.Range("I122").Value = 7040.875
$ws.Range("J122").Value = 3858.6
$ws.Range("K122").Value = 21122.625
$ws.Range("L122").Value = 11575.8
$ws.Range("M122").Value = -18672.625
$ws.Range("N122").Value = -16475.8

$ws.Range("H136").Value = 4851.816
$ws.Range("I136").Value = 1403.2858
$ws.Range("J136").Value = 9111.764999999999
$ws.Range("K136").Value = 4209.857400000001
$ws.Range("L136").Value = 27335.295
$ws.Range("M136").Value = -1659.857400000001
$ws.Range("N136").Value = -32435.295
